$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c_D2 = $ws.Range("D2")
$c_D2.NumberFormat = "@"
$c_D2.Value = "59.155.45"
$c_D2.Style = "Normal"

$c_E2 = $ws.Range("E2")
$c_E2.NumberFormat = "@"
$c_E2.Value = "  +3.64%  "
$c_E2.Style = "Normal"

$c_D3 = $ws.Range("D3")
$c_D3.NumberFormat = "@"
$c_D3.Value = "2.990.12"
$c_D3.Style = "Normal"

$c_E3 = $ws.Range("E3")
$c_E3.NumberFormat = "@"
$c_E3.Value = "  +3.32%  "
$c_E3.Style = "Normal"

$c_E4 = $ws.Range("E4")
$c_E4.NumberFormat = "@"
$c_E4.Value = "  -0.06%  "
$c_E4.Style = "Normal"

$c_D5 = $ws.Range("D5")
$c_D5.NumberFormat = "@"
$c_D5.Value = "562.81"
$c_D5.Style = "Normal"

$c_E5 = $ws.Range("E5")
$c_E5.NumberFormat = "@"
$c_E5.Value = "  +2.96%  "
$c_E5.Style = "Normal"

$c_D6 = $ws.Range("D6")
$c_D6.NumberFormat = "@"
$c_D6.Value = "138.27"
$c_D6.Style = "Normal"

$c_E6 = $ws.Range("E6")
$c_E6.NumberFormat = "@"
$c_E6.Value = "  +11.31%  "
$c_E6.Style = "Normal"

$c_E7 = $ws.Range("E7")
$c_E7.NumberFormat = "@"
$c_E7.Value = "  -0.19%  "
$c_E7.Style = "Normal"

$c_D8 = $ws.Range("D8")
$c_D8.NumberFormat = "@"
$c_D8.Value = "0.520"
$c_D8.Style = "Normal"

$c_E8 = $ws.Range("E8")
$c_E8.NumberFormat = "@"
$c_E8.Value = "  +4.10%  "
$c_E8.Style = "Normal"

$c_D9 = $ws.Range("D9")
$c_D9.NumberFormat = "@"
$c_D9.Value = "2.981.91"
$c_D9.Style = "Normal"

$c_E9 = $ws.Range("E9")
$c_E9.NumberFormat = "@"
$c_E9.Value = "  +3.20%  "
$c_E9.Style = "Normal"

$c_E10 = $ws.Range("E10")
$c_E10.NumberFormat = "@"
$c_E10.Value = "  +8.59%  "
$c_E10.Style = "Normal"

$c_D11 = $ws.Range("D11")
$c_D11.NumberFormat = "@"
$c_D11.Value = "5.12"
$c_D11.Style = "Normal"

$c_E11 = $ws.Range("E11")
$c_E11.NumberFormat = "@"
$c_E11.Value = "  +10.30%  "
$c_E11.Style = "Normal"

$c_E12 = $ws.Range("E12")
$c_E12.NumberFormat = "@"
$c_E12.Value = "  +5.36%  "
$c_E12.Style = "Normal"

$c_D13 = $ws.Range("D13")
$c_D13.NumberFormat = "@"
$c_D13.Value = "0.0000229"
$c_D13.Style = "Normal"

$c_E13 = $ws.Range("E13")
$c_E13.NumberFormat = "@"
$c_E13.Value = "  +9.85%  "
$c_E13.Style = "Normal"

$c_D14 = $ws.Range("D14")
$c_D14.NumberFormat = "@"
$c_D14.Value = "33.70"
$c_D14.Style = "Normal"

$c_E14 = $ws.Range("E14")
$c_E14.NumberFormat = "@"
$c_E14.Value = "  +5.26%  "
$c_E14.Style = "Normal"

$c_E15 = $ws.Range("E15")
$c_E15.NumberFormat = "@"
$c_E15.Value = "  +2.85%  "
$c_E15.Style = "Normal"

$c_D16 = $ws.Range("D16")
$c_D16.NumberFormat = "@"
$c_D16.Value = "3.490.41"
$c_D16.Style = "Normal"

$c_E16 = $ws.Range("E16")
$c_E16.NumberFormat = "@"
$c_E16.Value = "  +3.24%  "
$c_E16.Style = "Normal"

$c_D17 = $ws.Range("D17")
$c_D17.NumberFormat = "@"
$c_D17.Value = "7.04"
$c_D17.Style = "Normal"

$c_E17 = $ws.Range("E17")
$c_E17.NumberFormat = "@"
$c_E17.Value = "  +8.87%  "
$c_E17.Style = "Normal"

$c_D18 = $ws.Range("D18")
$c_D18.NumberFormat = "@"
$c_D18.Value = "2.992.45"
$c_D18.Style = "Normal"

$c_E18 = $ws.Range("E18")
$c_E18.NumberFormat = "@"
$c_E18.Value = "  +3.16%  "
$c_E18.Style = "Normal"

$c_D19 = $ws.Range("D19")
$c_D19.NumberFormat = "@"
$c_D19.Value = "59.121.71"
$c_D19.Style = "Normal"

$c_E19 = $ws.Range("E19")
$c_E19.NumberFormat = "@"
$c_E19.Value = "  +3.41%  "
$c_E19.Style = "Normal"

$c_D20 = $ws.Range("D20")
$c_D20.NumberFormat = "@"
$c_D20.Value = "427.90"
$c_D20.Style = "Normal"

$c_E20 = $ws.Range("E20")
$c_E20.NumberFormat = "@"
$c_E20.Value = "  +6.28%  "
$c_E20.Style = "Normal"

$c_D21 = $ws.Range("D21")
$c_D21.NumberFormat = "@"
$c_D21.Value = "13.54"
$c_D21.Style = "Normal"

$c_E21 = $ws.Range("E21")
$c_E21.NumberFormat = "@"
$c_E21.Value = "  +6.38%  "
$c_E21.Style = "Normal"

$c_E22 = $ws.Range("E22")
$c_E22.NumberFormat = "@"
$c_E22.Value = "  +7.21%  "
$c_E22.Style = "Normal"

$c_D23 = $ws.Range("D23")
$c_D23.NumberFormat = "@"
$c_D23.Value = "13.44"
$c_D23.Style = "Normal"

$c_E23 = $ws.Range("E23")
$c_E23.NumberFormat = "@"
$c_E23.Value = "  +7.16%  "
$c_E23.Style = "Normal"

$c_D24 = $ws.Range("D24")
$c_D24.NumberFormat = "@"
$c_D24.Value = "7.08"
$c_D24.Style = "Normal"

$c_E24 = $ws.Range("E24")
$c_E24.NumberFormat = "@"
$c_E24.Value = "  +4.38%  "
$c_E24.Style = "Normal"

$c_D25 = $ws.Range("D25")
$c_D25.NumberFormat = "@"
$c_D25.Value = "80.51"
$c_D25.Style = "Normal"

$c_E25 = $ws.Range("E25")
$c_E25.NumberFormat = "@"
$c_E25.Value = "  +4.20%  "
$c_E25.Style = "Normal"

$c_E26 = $ws.Range("E26")
$c_E26.NumberFormat = "@"
$c_E26.Value = "  +0.03%  "
$c_E26.Style = "Normal"

$c_E27 = $ws.Range("E27")
$c_E27.NumberFormat = "@"
$c_E27.Value = "  +0.00%  "
$c_E27.Style = "Normal"

$c_D28 = $ws.Range("D28")
$c_D28.NumberFormat = "@"
$c_D28.Value = "2.14"
$c_D28.Style = "Normal"

$c_E28 = $ws.Range("E28")
$c_E28.NumberFormat = "@"
$c_E28.Value = "  +12.23%  "
$c_E28.Style = "Normal"

$c_E29 = $ws.Range("E29")
$c_E29.NumberFormat = "@"
$c_E29.Value = "  +3.77%  "
$c_E29.Style = "Normal"

$c_E30 = $ws.Range("E30")
$c_E30.NumberFormat = "@"
$c_E30.Value = "  +7.80%  "
$c_E30.Style = "Normal"

$c_E31 = $ws.Range("E31")
$c_E31.NumberFormat = "@"
$c_E31.Value = "  +4.56%  "
$c_E31.Style = "Normal"

$c_D32 = $ws.Range("D32")
$c_D32.NumberFormat = "@"
$c_D32.Value = "6.12"
$c_D32.Style = "Normal"

$c_E32 = $ws.Range("E32")
$c_E32.NumberFormat = "@"
$c_E32.Value = "  +3.85%  "
$c_E32.Style = "Normal"

$c_D33 = $ws.Range("D33")
$c_D33.NumberFormat = "@"
$c_D33.Value = "0.0982"
$c_D33.Style = "Normal"

$c_E33 = $ws.Range("E33")
$c_E33.NumberFormat = "@"
$c_E33.Value = "  +0.40%  "
$c_E33.Style = "Normal"

$c_E34 = $ws.Range("E34")
$c_E34.NumberFormat = "@"
$c_E34.Value = "  +9.44%  "
$c_E34.Style = "Normal"

$c_D35 = $ws.Range("D35")
$c_D35.NumberFormat = "@"
$c_D35.Value = "0.0₃0770"
$c_D35.Style = "Normal"

$c_E35 = $ws.Range("E35")
$c_E35.NumberFormat = "@"
$c_E35.Value = "  +23.85%  "
$c_E35.Style = "Normal"

$c_E36 = $ws.Range("E36")
$c_E36.NumberFormat = "@"
$c_E36.Value = "  +7.59%  "
$c_E36.Style = "Normal"

$c_E37 = $ws.Range("E37")
$c_E37.NumberFormat = "@"
$c_E37.Value = "  +3.83%  "
$c_E37.Style = "Normal"

$c_D38 = $ws.Range("D38")
$c_D38.NumberFormat = "@"
$c_D38.Value = "49.02"
$c_D38.Style = "Normal"

$c_E38 = $ws.Range("E38")
$c_E38.NumberFormat = "@"
$c_E38.Value = "  +2.09%  "
$c_E38.Style = "Normal"

$c_E39 = $ws.Range("E39")
$c_E39.NumberFormat = "@"
$c_E39.Value = "  +5.93%  "
$c_E39.Style = "Normal"

$c_D40 = $ws.Range("D40")
$c_D40.NumberFormat = "@"
$c_D40.Value = "2.72"
$c_D40.Style = "Normal"

$c_E40 = $ws.Range("E40")
$c_E40.NumberFormat = "@"
$c_E40.Value = "  +13.42%  "
$c_E40.Style = "Normal"

$c_D41 = $ws.Range("D41")
$c_D41.NumberFormat = "@"
$c_D41.Value = "398.13"
$c_D41.Style = "Normal"

$c_E41 = $ws.Range("E41")
$c_E41.NumberFormat = "@"
$c_E41.Value = "  +10.51%  "
$c_E41.Style = "Normal"

$c_E42 = $ws.Range("E42")
$c_E42.NumberFormat = "@"
$c_E42.Value = "  +4.65%  "
$c_E42.Style = "Normal"

$c_D43 = $ws.Range("D43")
$c_D43.NumberFormat = "@"
$c_D43.Value = "2.748.40"
$c_D43.Style = "Normal"

$c_E43 = $ws.Range("E43")
$c_E43.NumberFormat = "@"
$c_E43.Value = "  +5.04%  "
$c_E43.Style = "Normal"

$c_E44 = $ws.Range("E44")
$c_E44.NumberFormat = "@"
$c_E44.Value = "  +1.55%  "
$c_E44.Style = "Normal"

$c_D45 = $ws.Range("D45")
$c_D45.NumberFormat = "@"
$c_D45.Value = "0.252"
$c_D45.Style = "Normal"

$c_E45 = $ws.Range("E45")
$c_E45.NumberFormat = "@"
$c_E45.Value = "  +11.23%  "
$c_E45.Style = "Normal"

$c_D47 = $ws.Range("D47")
$c_D47.NumberFormat = "@"
$c_D47.Value = "122.66"
$c_D47.Style = "Normal"

$c_E47 = $ws.Range("E47")
$c_E47.NumberFormat = "@"
$c_E47.Value = "  +3.05%  "
$c_E47.Style = "Normal"

$c_E48 = $ws.Range("E48")
$c_E48.NumberFormat = "@"
$c_E48.Value = "  +2.83%  "
$c_E48.Style = "Normal"

$c_E49 = $ws.Range("E49")
$c_E49.NumberFormat = "@"
$c_E49.Value = "  +4.38%  "
$c_E49.Style = "Normal"

$c_B50 = $ws.Range("B50")
$c_B50.NumberFormat = "@"
$c_B50.Value = "Arweave"
$c_B50.Style = "Normal"

$c_C50 = $ws.Range("C50")
$c_C50.NumberFormat = "@"
$c_C50.Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$c_C50.Style = "Normal"

$c_D50 = $ws.Range("D50")
$c_D50.NumberFormat = "@"
$c_D50.Value = "32.35"
$c_D50.Style = "Normal"

$c_E50 = $ws.Range("E50")
$c_E50.NumberFormat = "@"
$c_E50.Value = "  +20.05%  "
$c_E50.Style = "Normal"

$c_B51 = $ws.Range("B51")
$c_B51.NumberFormat = "@"
$c_B51.Value = "InjectiveProtocol"
$c_B51.Style = "Normal"

$c_C51 = $ws.Range("C51")
$c_C51.NumberFormat = "@"
$c_C51.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c_C51.Style = "Normal"

$c_D51 = $ws.Range("D51")
$c_D51.NumberFormat = "@"
$c_D51.Value = "23.45"
$c_D51.Style = "Normal"

$c_E51 = $ws.Range("E51")
$c_E51.NumberFormat = "@"
$c_E51.Value = "  +4.71%  "
$c_E51.Style = "Normal"

